$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the 8 new rows (120-127) describing new "Salvador Allende" place entries ---
$ws.Range("A120").Value2 = 119
$ws.Range("B120").Value2 = "Salvador Allende street"
$ws.Range("C120").Value2 = "street"
$ws.Range("D120").Value2 = "Europe"
$ws.Range("E120").Value2 = "Bosnia and Herzegovina"
$ws.Range("F120").Value2 = "Sarajevo"
$ws.Range("G120").Value2 = "Dobrinja"
$ws.Range("N120").Value2 = 2007
$ws.Range("O120").Value2 = 2
$ws.Range("P120").Value2 = 18
$ws.Range("V120").Value2 = 0
$ws.Range("Y120").Value2 = "http://www.abacq.org/calle/index.php?2007/02/18/59-sarajevo-bosnia-herzegovina"
$ws.Range("A121").Value2 = 120
$ws.Range("B121").Value2 = "Villa Salvador Allende"
$ws.Range("C121").Value2 = "neighborhood"
$ws.Range("D121").Value2 = "South America"
$ws.Range("E121").Value2 = "Chile"
$ws.Range("F121").Value2 = "Región Metropolitana"
$ws.Range("G121").Value2 = "Santiago"
$ws.Range("H121").Value2 = "La Pintana"
$ws.Range("K121").Value2 = 8010000
$ws.Range("L121").Value2 = -33.580694899999997
$ws.Range("M121").Value2 = -70.646144500000005
$ws.Range("N121").Value2 = 2012
$ws.Range("O121").Value2 = 11
$ws.Range("Q121").Value2 = "google maps"
$ws.Range("V121").Value2 = 1
$ws.Range("W121").Value2 = "https://www.openstreetmap.org/node/9728569027"
$ws.Range("X121").Value2 = "https://goo.gl/maps/NbqWJV69Cbu8h53k6"
$ws.Range("A122").Value2 = 121
$ws.Range("B122").Value2 = "Praça Salvador Allende"
$ws.Range("C122").Value2 = "park"
$ws.Range("D122").Value2 = "South America"
$ws.Range("E122").Value2 = "Brazil"
$ws.Range("F122").Value2 = "Rio Grande do Sul"
$ws.Range("G122").Value2 = "Porto Alegre"
$ws.Range("H122").Value2 = "Centro Histórico"
$ws.Range("K122").Value2 = "90050-200"
$ws.Range("L122").Value2 = -30.034959675270098
$ws.Range("M122").Value2 = -51.223894258388903
$ws.Range("N122").Value2 = 2004
$ws.Range("O122").Value2 = 9
$ws.Range("Q122").Value2 = "desc abacq"
$ws.Range("R122").Value2 = "Prefeitura Municipal de Porto Alegre`nSecretaria Municipal do Meio Ambiente`nPraça Salvador Allende `nPresidente chileno. Herói das Américas`n2013"
$ws.Range("S122").Value2 = "pt"
$ws.Range("V122").Value2 = 1
$ws.Range("W122").Value2 = "https://www.openstreetmap.org/way/305752671"
$ws.Range("X122").Value2 = "https://goo.gl/maps/4FGZmC5aXyT4RUSc9"
$ws.Range("Y122").Value2 = "http://www.abacq.org/calle/index.php?2013/06/24/597-porto-alegre-brasil"
$ws.Range("A123").Value2 = 122
$ws.Range("B123").Value2 = "Monumento al presidente Salvador Allende en Brasilia"
$ws.Range("C123").Value2 = "monument"
$ws.Range("D123").Value2 = "South America"
$ws.Range("E123").Value2 = "Brazil"
$ws.Range("F123").Value2 = "Distrito Federal"
$ws.Range("G123").Value2 = "Brasilia"
$ws.Range("N123").Value2 = 2008
$ws.Range("Q123").Value2 = "desc abacq"
$ws.Range("R123").Value2 = "Salvador Allende`nTenho fé no Chile e no seu destino, sigam vocês sabendo que, muito mais cedo do que tarde, abrir-se-âo as grandes alamedas por onde passe o homem livre, para construir uma sociedade melhor.`n(Presidente Allende, 11 de setembro 1973)"
$ws.Range("S123").Value2 = "pt"
$ws.Range("V123").Value2 = 0
$ws.Range("Y123").Value2 = "http://www.abacq.org/calle/index.php?2011/12/12/550-brasilia-brasil"
$ws.Range("A124").Value2 = 123
$ws.Range("B124").Value2 = "Avenida Salvador Allende"
$ws.Range("C124").Value2 = "street"
$ws.Range("D124").Value2 = "South America"
$ws.Range("E124").Value2 = "Brazil"
$ws.Range("F124").Value2 = "Rio de Janeiro"
$ws.Range("G124").Value2 = "Rio de Janeiro"
$ws.Range("K124").Value2 = "22783-020"
$ws.Range("L124").Value2 = -22.985492037936801
$ws.Range("M124").Value2 = -43.413389967982098
$ws.Range("N124").Value2 = 2007
$ws.Range("O124").Value2 = 2
$ws.Range("P124").Value2 = 18
$ws.Range("Q124").Value2 = "abacq date posted"
$ws.Range("R124").Value2 = "Avenida Salvador Allende`n(1908 - 1973) Médico, fundador do Partido Socialista, presidente do Chile (1970-1973)."
$ws.Range("S124").Value2 = "pt"
$ws.Range("V124").Value2 = 1
$ws.Range("W124").Value2 = "https://www.openstreetmap.org/way/426694183"
$ws.Range("X124").Value2 = "https://goo.gl/maps/fakGRb84JRo2jQj8A"
$ws.Range("Y124").Value2 = "http://www.abacq.org/calle/index.php?2010/07/29/499-rio-de-janeiro-brasil"
$ws.Range("A125").Value2 = 124
$ws.Range("B125").Value2 = "Praça Salvador Allende"
$ws.Range("C125").Value2 = "street"
$ws.Range("D125").Value2 = "South America"
$ws.Range("E125").Value2 = "Brazil"
$ws.Range("F125").Value2 = "São Paulo"
$ws.Range("G125").Value2 = "São Paulo"
$ws.Range("H125").Value2 = "Sacomã"
$ws.Range("I125").Value2 = "Vila das Mercês"
$ws.Range("K125").Value2 = "04165-000"
$ws.Range("L125").Value2 = -23.618749999999999
$ws.Range("M125").Value2 = -46.610177
$ws.Range("N125").Value2 = 2007
$ws.Range("O125").Value2 = 2
$ws.Range("P125").Value2 = 18
$ws.Range("Q125").Value2 = "abacq date posted"
$ws.Range("V125").Value2 = 1
$ws.Range("W125").Value2 = "https://www.openstreetmap.org/way/292191474"
$ws.Range("X125").Value2 = "https://goo.gl/maps/oCKr4cgZmgq7mkCQ6"
$ws.Range("Y125").Value2 = "http://www.abacq.org/calle/index.php?2008/12/08/291-sao-paulo-brasil"
$ws.Range("A126").Value2 = 125
$ws.Range("B126").Value2 = "Salvador Allende Bus Station"
$ws.Range("C126").Value2 = "public transport station"
$ws.Range("D126").Value2 = "South America"
$ws.Range("E126").Value2 = "Brazil"
$ws.Range("F126").Value2 = "Rio de Janeiro"
$ws.Range("G126").Value2 = "Rio de Janeiro"
$ws.Range("H126").Value2 = "Recreio dos Bandeirantes"
$ws.Range("I126").Value2 = "Avenida das Américas"
$ws.Range("K126").Value2 = "22790-710"
$ws.Range("L126").Value2 = -23.0082369008645
$ws.Range("M126").Value2 = -43.442585473568698
$ws.Range("N126").Value2 = 2014
$ws.Range("O126").Value2 = 4
$ws.Range("P126").Value2 = 22
$ws.Range("Q126").Value2 = "openstreetmap"
$ws.Range("V126").Value2 = 1
$ws.Range("W126").Value2 = "https://www.openstreetmap.org/way/276421960"
$ws.Range("X126").Value2 = "https://goo.gl/maps/hfjxEk8grbHLu8239"
$ws.Range("A127").Value2 = 126
$ws.Range("B127").Value2 = "бул. „Андрей Сахаров“ / Boulevard `"Andrej Sakharov`""
$ws.Range("C127").Value2 = "street"
$ws.Range("D127").Value2 = "Europe"
$ws.Range("E127").Value2 = "Bulgaria"
$ws.Range("F127").Value2 = "Sofia City Province"
$ws.Range("G127").Value2 = "Sofia"
$ws.Range("H127").Value2 = "Mladost"
$ws.Range("I127").Value2 = "Mladost 1"
$ws.Range("K127").Value2 = 1729
$ws.Range("L127").Value2 = 42.652000000000001
$ws.Range("M127").Value2 = 23.37491
$ws.Range("N127").Value2 = 1977
$ws.Range("Q127").Value2 = "desc abacq"
$ws.Range("R127").Value2 = "Salvador Allende, 1908-1973, hijo ilustre del pueblo chileno`nEntregó su vida por la causa de la paz, la democracia y el progreso socialista"
$ws.Range("S127").Value2 = "bg"
$ws.Range("U127").Value2 = "бул. „Салвадор Алиенде“ / Boulevard `"Salvador Allende`""
$ws.Range("V127").Value2 = 1
$ws.Range("W127").Value2 = "https://www.openstreetmap.org/way/69596345"
$ws.Range("X127").Value2 = "https://goo.gl/maps/LfxjeoGtcJmf8gm16"
$ws.Range("Y127").Value2 = "http://www.abacq.org/calle/index.php?2007/10/02/110-sofia-bulgaria"

# --- Refresh the AutoFilter range to cover the newly added rows ---
$ws.AutoFilterMode = $false
$ws.Range("B1:AO127").AutoFilter()

# --- Keep the _FilterDatabase defined name in sync with the new filter range ---
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=allende!`$B`$1:`$AO`$127"
    }
}

# --- Mirror the author's final selection (cell just past the last data row) ---
$ws.Range("A128").Select()
